$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-ValueCell($rowIndex, $value) {
    $row = $t.Rows.Item($rowIndex)
    $c = $row.Cells.Item(2)
    $r = $c.Range
    $r.Text = $value
    $r.Font.Bold = 1
    $r.Font.Size = 12
    $r.Font.SizeBi = 12
}

# RETENTION Ratio
Set-ValueCell 24 "0.8571"

# Answer Recall Lenient (ARL)
Set-ValueCell 44 "0.375"

# Answer Recall Strict (ARS)
Set-ValueCell 45 "0.25"

# Answer Recall Average (ARA)
Set-ValueCell 46 "0.3125"
